$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "LastName" / "FirstName" column headers
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

# Renumber the CNE column (A2:A11): 19000031-19000040 -> 19000001-19000010
for ($i = 2; $i -le 11; $i++) {
    $ws.Cells.Item($i, 1).Value = 19000000 + ($i - 1)
}

# Reset the header row back to the plain/default look, then (re)stamp a
# fresh cell style across the header + id/name columns (A1:C11) so they all
# share one explicit style separate from the untouched columns to the right
$ws.Range("A1:C1").ClearFormats()
$ws.Range("A1:C11").VerticalAlignment = -4107

# Update the view: scroll back to A1 and move the selection to F8
$ws.Range("F8").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
